# Fruta / hortaliza, semanal
# Insert 4 new weekly-report rows (Plátano, Vega Central Mapocho de Santiago)
# just above the existing row 600, pushing the remaining rows down
# (old rows 600-693 become 604-697).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 600..603.
$ws.Range("A600:A603").EntireRow.Insert()

# Shared/boilerplate values for this market+product block (same for every
# row in this sheet).
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$prodId    = 100108
$producto  = "Tropicales y subtropicales"
$catId     = 100108006
$categoria = "Plátano"
$unidad    = "`$/caja 20 kilos"
$origen    = "Ecuador"
$kgUnidad  = 20

# Row 600
$ws.Cells.Item(600, 1).Value  = $mercadoId
$ws.Cells.Item(600, 2).Value  = $mercado
$ws.Cells.Item(600, 3).Value  = $region
$ws.Cells.Item(600, 4).Value  = 44474
$ws.Cells.Item(600, 5).Value  = $codreg
$ws.Cells.Item(600, 6).Value  = $tipo
$ws.Cells.Item(600, 7).Value  = $prodId
$ws.Cells.Item(600, 8).Value  = $producto
$ws.Cells.Item(600, 9).Value  = $catId
$ws.Cells.Item(600, 10).Value = $categoria
$ws.Cells.Item(600, 11).Value = "Sin especificar"
$ws.Cells.Item(600, 12).Value = "Maduro"
$ws.Cells.Item(600, 13).Value = 300
$ws.Cells.Item(600, 14).Value = 17000
$ws.Cells.Item(600, 15).Value = 17000
$ws.Cells.Item(600, 16).Value = 17000
$ws.Cells.Item(600, 17).Value = $unidad
$ws.Cells.Item(600, 18).Value = $origen
$ws.Cells.Item(600, 19).Value = 850
$ws.Cells.Item(600, 20).Value = $kgUnidad

# Row 601
$ws.Cells.Item(601, 1).Value  = $mercadoId
$ws.Cells.Item(601, 2).Value  = $mercado
$ws.Cells.Item(601, 3).Value  = $region
$ws.Cells.Item(601, 4).Value  = 44474
$ws.Cells.Item(601, 5).Value  = $codreg
$ws.Cells.Item(601, 6).Value  = $tipo
$ws.Cells.Item(601, 7).Value  = $prodId
$ws.Cells.Item(601, 8).Value  = $producto
$ws.Cells.Item(601, 9).Value  = $catId
$ws.Cells.Item(601, 10).Value = $categoria
$ws.Cells.Item(601, 11).Value = "Sin especificar"
$ws.Cells.Item(601, 12).Value = "Pintón"
$ws.Cells.Item(601, 13).Value = 450
$ws.Cells.Item(601, 14).Value = 18000
$ws.Cells.Item(601, 15).Value = 18000
$ws.Cells.Item(601, 16).Value = 18000
$ws.Cells.Item(601, 17).Value = $unidad
$ws.Cells.Item(601, 18).Value = $origen
$ws.Cells.Item(601, 19).Value = 900
$ws.Cells.Item(601, 20).Value = $kgUnidad

# Row 602
$ws.Cells.Item(602, 1).Value  = $mercadoId
$ws.Cells.Item(602, 2).Value  = $mercado
$ws.Cells.Item(602, 3).Value  = $region
$ws.Cells.Item(602, 4).Value  = 44474
$ws.Cells.Item(602, 5).Value  = $codreg
$ws.Cells.Item(602, 6).Value  = $tipo
$ws.Cells.Item(602, 7).Value  = $prodId
$ws.Cells.Item(602, 8).Value  = $producto
$ws.Cells.Item(602, 9).Value  = $catId
$ws.Cells.Item(602, 10).Value = $categoria
$ws.Cells.Item(602, 11).Value = "Sin especificar"
$ws.Cells.Item(602, 12).Value = "Primera Maduro"
$ws.Cells.Item(602, 13).Value = 400
$ws.Cells.Item(602, 14).Value = 19000
$ws.Cells.Item(602, 15).Value = 19000
$ws.Cells.Item(602, 16).Value = 19000
$ws.Cells.Item(602, 17).Value = $unidad
$ws.Cells.Item(602, 18).Value = $origen
$ws.Cells.Item(602, 19).Value = 950
$ws.Cells.Item(602, 20).Value = $kgUnidad

# Row 603
$ws.Cells.Item(603, 1).Value  = $mercadoId
$ws.Cells.Item(603, 2).Value  = $mercado
$ws.Cells.Item(603, 3).Value  = $region
$ws.Cells.Item(603, 4).Value  = 44474
$ws.Cells.Item(603, 5).Value  = $codreg
$ws.Cells.Item(603, 6).Value  = $tipo
$ws.Cells.Item(603, 7).Value  = $prodId
$ws.Cells.Item(603, 8).Value  = $producto
$ws.Cells.Item(603, 9).Value  = $catId
$ws.Cells.Item(603, 10).Value = $categoria
$ws.Cells.Item(603, 11).Value = "Sin especificar"
$ws.Cells.Item(603, 12).Value = "Primera Pintón"
$ws.Cells.Item(603, 13).Value = 350
$ws.Cells.Item(603, 14).Value = 20000
$ws.Cells.Item(603, 15).Value = 20000
$ws.Cells.Item(603, 16).Value = 20000
$ws.Cells.Item(603, 17).Value = $unidad
$ws.Cells.Item(603, 18).Value = $origen
$ws.Cells.Item(603, 19).Value = 1000
$ws.Cells.Item(603, 20).Value = $kgUnidad
